$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, shifting existing rows 91-102 down to 92-103.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new data record.
$ws.Cells.Item(91, 1).Value = 9
$ws.Cells.Item(91, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(91, 3).Value = "Metropolitana"
$ws.Cells.Item(91, 4).Value = 45127
$ws.Cells.Item(91, 5).Value = 13
$ws.Cells.Item(91, 6).Value = 100112035
$ws.Cells.Item(91, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 26
$ws.Cells.Item(91, 11).Value = 17000
$ws.Cells.Item(91, 12).Value = 17000
$ws.Cells.Item(91, 13).Value = 17000
$ws.Cells.Item(91, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(91, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(91, 16).Value = 1133
$ws.Cells.Item(91, 17).Value = 15
$ws.Cells.Item(91, 18).Value = "Hortaliza"
